$wb = $excel.ActiveWorkbook

# --- Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("M2").Value = 0.01300380348513583
$ws.Range("N2").Value = 0.06470716043000752
$ws.Range("O2").Value = 0.1768920421031213
$ws.Range("Q2").Value = 0.05352418943867731
$ws.Range("R2").Value = 0.1025440908090517
$ws.Range("S2").Value = 0.0302575094035128
$ws.Range("U2").Value = 0.1168701292929945
$ws.Range("V2").Value = 0.2060485475538608
$ws.Range("X2").Value = 0.05785250123565147
$ws.Range("Z2").Value = 0.02963409442633952
$ws.Range("AA2").Value = 0.001125721343911373
$ws.Range("AB2").Value = 0.01243689992976523
$ws.Range("AC2").Value = 0.01682819365710547
$ws.Range("AE2").Value = 0.02826435346128461
$ws.Range("AF2").Value = 0.03491450398798339
$ws.Range("AH2").Value = 0.008973415952351228
$ws.Range("AJ2").Value = 0.007073109023518038
$ws.Range("AK2").Value = 0.01563887754088136
$ws.Range("AL2").Value = 0.01395056229331972
$ws.Range("AM2").Value = 0.004893159014755649
$ws.Range("AP2").Value = 0.004567135616771063
$ws.Range("E3").Value = 0.2741814000016149
$ws.Range("F3").Value = 0.07443462849669671
$ws.Range("G3").Value = 0.1509708869539654
$ws.Range("I3").Value = 0.1472449326097034
$ws.Range("K3").Value = 0.02094875236179295
$ws.Range("L3").Value = 0.07368266161291535
$ws.Range("M3").Value = 0.149209551886532
$ws.Range("N3").Value = 0.02392896891116136
$ws.Range("Q3").Value = 0.009294198521904442
$ws.Range("T3").Value = 0.001420849088717846
$ws.Range("W3").Value = 0.0154934716650069
$ws.Range("X3").Value = 0.01122546621864679
$ws.Range("AB3").Value = 0.02565361780472893
$ws.Range("AC3").Value = 0.01003981970945778
$ws.Range("AE3").Value = 0.006946288660628832
$ws.Range("AH3").Value = 0.00532450549652647
$ws.Range("E4").Value = 0.1854996753135776
$ws.Range("F4").Value = 0.02137935804147948
$ws.Range("G4").Value = 0.1562103698850893
$ws.Range("H4").Value = 0.02273371733624137
$ws.Range("I4").Value = 0.09922162740270939
$ws.Range("K4").Value = 0.01492043414513093
$ws.Range("L4").Value = 0.03360976426005838
$ws.Range("M4").Value = 0.1687125481789583
$ws.Range("N4").Value = 0.05013855319846904
$ws.Range("Q4").Value = 0.005237391215597353
$ws.Range("R4").Value = 0.006620375154773804
$ws.Range("U4").Value = 0.02914737337619114
$ws.Range("W4").Value = 0.01828068700667387
$ws.Range("X4").Value = 0.04932887488570712
$ws.Range("AB4").Value = 0.07660755328967615
$ws.Range("AC4").Value = 0.04776304182488846
$ws.Range("AE4").Value = 0.01429122955439394
$ws.Range("AJ4").Value = 0.0002974259303844077
$ws.Range("M5").Value = 0.1619070324096959
$ws.Range("O5").Value = 0.206635617414818
$ws.Range("Q5").Value = 0.07664663686483485
$ws.Range("R5").Value = 0.02791813557654366
$ws.Range("S5").Value = 0.01240428936138394
$ws.Range("T5").Value = 0.07024392445133758
$ws.Range("U5").Value = 0.215770792750347
$ws.Range("V5").Value = 0.06321525403400363
$ws.Range("Y5").Value = 0.04263734784852644
$ws.Range("AC5").Value = 0.01864208134227395
$ws.Range("AE5").Value = 0.0289071488779056
$ws.Range("AF5").Value = 0.03726537467397994
$ws.Range("AK5").Value = 0.02429090895649751
$ws.Range("AM5").Value = 0.0005361060469142327
$ws.Range("AP5").Value = 0.01297934939093773
$ws.Range("D6").Value = 0.0187109963177049
$ws.Range("E6").Value = 0.1177015998425562
$ws.Range("F6").Value = 0.276455136601797
$ws.Range("G6").Value = 0.03638386751704333
$ws.Range("H6").Value = 0.04698146537710059
$ws.Range("I6").Value = 0.08152057879481077
$ws.Range("K6").Value = 0.01778441426110976
$ws.Range("L6").Value = 0.1365668852907765
$ws.Range("M6").Value = 0.09107721257801897
$ws.Range("P6").Value = 0.003513404858519884
$ws.Range("Q6").Value = 0.007940276467537252
$ws.Range("T6").Value = 0.02289006985968895
$ws.Range("V6").Value = 0.003386220553015345
$ws.Range("W6").Value = 0.04983653427890315
$ws.Range("AA6").Value = 0.005713916020895102
$ws.Range("AB6").Value = 0.06597176324457092
$ws.Range("AC6").Value = 0.002410984200145435
$ws.Range("AE6").Value = 0.004728449972345301
$ws.Range("AG6").Value = 0.00858061769899039
$ws.Range("AH6").Value = 0.00184560626447008

# --- Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("M2").Value = 0.01300380348513583
$ws.Range("N2").Value = 0.07771096391514334
$ws.Range("O2").Value = 0.2546030060182646
$ws.Range("P2").Value = 0.2546030060182646
$ws.Range("Q2").Value = 0.3081271954569419
$ws.Range("R2").Value = 0.4106712862659936
$ws.Range("S2").Value = 0.4409287956695064
$ws.Range("T2").Value = 0.4409287956695064
$ws.Range("U2").Value = 0.5577989249625009
$ws.Range("V2").Value = 0.7638474725163618
$ws.Range("W2").Value = 0.7638474725163618
$ws.Range("X2").Value = 0.8216999737520132
$ws.Range("Y2").Value = 0.8216999737520132
$ws.Range("Z2").Value = 0.8513340681783527
$ws.Range("AA2").Value = 0.8524597895222641
$ws.Range("AB2").Value = 0.8648966894520294
$ws.Range("AC2").Value = 0.8817248831091348
$ws.Range("AD2").Value = 0.8817248831091348
$ws.Range("AE2").Value = 0.9099892365704194
$ws.Range("AF2").Value = 0.9449037405584028
$ws.Range("AG2").Value = 0.9449037405584028
$ws.Range("AH2").Value = 0.953877156510754
$ws.Range("AI2").Value = 0.953877156510754
$ws.Range("AJ2").Value = 0.960950265534272
$ws.Range("AK2").Value = 0.9765891430751533
$ws.Range("AL2").Value = 0.990539705368473
$ws.Range("AM2").Value = 0.9954328643832286
$ws.Range("AN2").Value = 0.9954328643832286
$ws.Range("AO2").Value = 0.9954328643832286
$ws.Range("AP2").Value = 0.9999999999999997
$ws.Range("AQ2").Value = 0.9999999999999997
$ws.Range("AR2").Value = 0.9999999999999997
$ws.Range("E3").Value = 0.2741814000016149
$ws.Range("F3").Value = 0.3486160284983116
$ws.Range("G3").Value = 0.4995869154522771
$ws.Range("H3").Value = 0.4995869154522771
$ws.Range("I3").Value = 0.6468318480619805
$ws.Range("J3").Value = 0.6468318480619805
$ws.Range("K3").Value = 0.6677806004237734
$ws.Range("L3").Value = 0.7414632620366888
$ws.Range("M3").Value = 0.8906728139232208
$ws.Range("N3").Value = 0.9146017828343822
$ws.Range("O3").Value = 0.9146017828343822
$ws.Range("P3").Value = 0.9146017828343822
$ws.Range("Q3").Value = 0.9238959813562867
$ws.Range("R3").Value = 0.9238959813562867
$ws.Range("S3").Value = 0.9238959813562867
$ws.Range("T3").Value = 0.9253168304450046
$ws.Range("U3").Value = 0.9253168304450046
$ws.Range("V3").Value = 0.9253168304450046
$ws.Range("W3").Value = 0.9408103021100115
$ws.Range("X3").Value = 0.9520357683286583
$ws.Range("Y3").Value = 0.9520357683286583
$ws.Range("Z3").Value = 0.9520357683286583
$ws.Range("AA3").Value = 0.9520357683286583
$ws.Range("AB3").Value = 0.9776893861333873
$ws.Range("AC3").Value = 0.9877292058428451
$ws.Range("AD3").Value = 0.9877292058428451
$ws.Range("AE3").Value = 0.9946754945034739
$ws.Range("AF3").Value = 0.9946754945034739
$ws.Range("AG3").Value = 0.9946754945034739
$ws.Range("E4").Value = 0.1854996753135776
$ws.Range("F4").Value = 0.2068790333550571
$ws.Range("G4").Value = 0.3630894032401464
$ws.Range("H4").Value = 0.3858231205763877
$ws.Range("I4").Value = 0.4850447479790971
$ws.Range("J4").Value = 0.4850447479790971
$ws.Range("K4").Value = 0.4999651821242281
$ws.Range("L4").Value = 0.5335749463842865
$ws.Range("M4").Value = 0.7022874945632447
$ws.Range("N4").Value = 0.7524260477617137
$ws.Range("O4").Value = 0.7524260477617137
$ws.Range("P4").Value = 0.7524260477617137
$ws.Range("Q4").Value = 0.757663438977311
$ws.Range("R4").Value = 0.7642838141320848
$ws.Range("S4").Value = 0.7642838141320848
$ws.Range("T4").Value = 0.7642838141320848
$ws.Range("U4").Value = 0.7934311875082759
$ws.Range("V4").Value = 0.7934311875082759
$ws.Range("W4").Value = 0.8117118745149497
$ws.Range("X4").Value = 0.8610407494006569
$ws.Range("Y4").Value = 0.8610407494006569
$ws.Range("Z4").Value = 0.8610407494006569
$ws.Range("AA4").Value = 0.8610407494006569
$ws.Range("AB4").Value = 0.9376483026903331
$ws.Range("AC4").Value = 0.9854113445152215
$ws.Range("AD4").Value = 0.9854113445152215
$ws.Range("AE4").Value = 0.9997025740696155
$ws.Range("AF4").Value = 0.9997025740696155
$ws.Range("AG4").Value = 0.9997025740696155
$ws.Range("AH4").Value = 0.9997025740696155
$ws.Range("AI4").Value = 0.9997025740696155
$ws.Range("AJ4").Value = 0.9999999999999999
$ws.Range("AK4").Value = 0.9999999999999999
$ws.Range("AL4").Value = 0.9999999999999999
$ws.Range("AM4").Value = 0.9999999999999999
$ws.Range("AN4").Value = 0.9999999999999999
$ws.Range("AO4").Value = 0.9999999999999999
$ws.Range("AP4").Value = 0.9999999999999999
$ws.Range("AQ4").Value = 0.9999999999999999
$ws.Range("AR4").Value = 0.9999999999999999
$ws.Range("M5").Value = 0.1619070324096959
$ws.Range("N5").Value = 0.1619070324096959
$ws.Range("O5").Value = 0.3685426498245139
$ws.Range("P5").Value = 0.3685426498245139
$ws.Range("Q5").Value = 0.4451892866893488
$ws.Range("R5").Value = 0.4731074222658924
$ws.Range("S5").Value = 0.4855117116272764
$ws.Range("T5").Value = 0.5557556360786139
$ws.Range("U5").Value = 0.7715264288289609
$ws.Range("V5").Value = 0.8347416828629646
$ws.Range("W5").Value = 0.8347416828629646
$ws.Range("X5").Value = 0.8347416828629646
$ws.Range("Y5").Value = 0.877379030711491
$ws.Range("Z5").Value = 0.877379030711491
$ws.Range("AA5").Value = 0.877379030711491
$ws.Range("AB5").Value = 0.877379030711491
$ws.Range("AC5").Value = 0.896021112053765
$ws.Range("AD5").Value = 0.896021112053765
$ws.Range("AE5").Value = 0.9249282609316706
$ws.Range("AF5").Value = 0.9621936356056505
$ws.Range("AG5").Value = 0.9621936356056505
$ws.Range("AH5").Value = 0.9621936356056505
$ws.Range("AI5").Value = 0.9621936356056505
$ws.Range("AJ5").Value = 0.9621936356056505
$ws.Range("AK5").Value = 0.9864845445621481
$ws.Range("AL5").Value = 0.9864845445621481
$ws.Range("AM5").Value = 0.9870206506090623
$ws.Range("AN5").Value = 0.9870206506090623
$ws.Range("AO5").Value = 0.9870206506090623
$ws.Range("D6").Value = 0.0187109963177049
$ws.Range("E6").Value = 0.1364125961602611
$ws.Range("F6").Value = 0.4128677327620581
$ws.Range("G6").Value = 0.4492516002791014
$ws.Range("H6").Value = 0.496233065656202
$ws.Range("I6").Value = 0.5777536444510127
$ws.Range("J6").Value = 0.5777536444510127
$ws.Range("K6").Value = 0.5955380587121225
$ws.Range("L6").Value = 0.732104944002899
$ws.Range("M6").Value = 0.823182156580918
$ws.Range("N6").Value = 0.823182156580918
$ws.Range("O6").Value = 0.823182156580918
$ws.Range("P6").Value = 0.8266955614394379
$ws.Range("Q6").Value = 0.8346358379069752
$ws.Range("R6").Value = 0.8346358379069752
$ws.Range("S6").Value = 0.8346358379069752
$ws.Range("T6").Value = 0.8575259077666642
$ws.Range("U6").Value = 0.8575259077666642
$ws.Range("V6").Value = 0.8609121283196794
$ws.Range("W6").Value = 0.9107486625985826
$ws.Range("X6").Value = 0.9107486625985826
$ws.Range("Y6").Value = 0.9107486625985826
$ws.Range("Z6").Value = 0.9107486625985826
$ws.Range("AA6").Value = 0.9164625786194777
$ws.Range("AB6").Value = 0.9824343418640485
$ws.Range("AC6").Value = 0.9848453260641939
$ws.Range("AD6").Value = 0.9848453260641939
$ws.Range("AE6").Value = 0.9895737760365393
$ws.Range("AF6").Value = 0.9895737760365393
$ws.Range("AG6").Value = 0.9981543937355297
$ws.Range("AH6").Value = 0.9999999999999998
$ws.Range("AI6").Value = 0.9999999999999998
$ws.Range("AJ6").Value = 0.9999999999999998
$ws.Range("AK6").Value = 0.9999999999999998
$ws.Range("AL6").Value = 0.9999999999999998
$ws.Range("AM6").Value = 0.9999999999999998
$ws.Range("AN6").Value = 0.9999999999999998
$ws.Range("AO6").Value = 0.9999999999999998
$ws.Range("AP6").Value = 0.9999999999999998
$ws.Range("AQ6").Value = 0.9999999999999998
$ws.Range("AR6").Value = 0.9999999999999998

# --- Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.5577989249625009
$ws.Range("D3").Value = 8
$ws.Range("F3").Value = 0.6468318480619805
$ws.Range("G3").Value = 6
$ws.Range("D4").Value = 11
$ws.Range("F4").Value = 0.5335749463842865
$ws.Range("G4").Value = 9
$ws.Range("D5").Value = 19
$ws.Range("F5").Value = 0.5557556360786139
$ws.Range("G5").Value = 9
$ws.Range("D6").Value = 8
$ws.Range("F6").Value = 0.5777536444510127
$ws.Range("G6").Value = 6

# --- Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.7638474725163618
$ws.Range("D3").Value = 11
$ws.Range("F3").Value = 0.7414632620366888
$ws.Range("G3").Value = 9
$ws.Range("F4").Value = 0.7022874945632447
$ws.Range("F5").Value = 0.7715264288289609
$ws.Range("F6").Value = 0.732104944002899

# --- Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 23
$ws.Range("F2").Value = 0.8216999737520132
$ws.Range("G2").Value = 12
$ws.Range("D3").Value = 12
$ws.Range("F3").Value = 0.8906728139232208
$ws.Range("G3").Value = 10
$ws.Range("D4").Value = 22
$ws.Range("F4").Value = 0.8117118745149497
$ws.Range("G4").Value = 20
$ws.Range("D5").Value = 21
$ws.Range("F5").Value = 0.8347416828629646
$ws.Range("G5").Value = 11
$ws.Range("D6").Value = 12
$ws.Range("F6").Value = 0.823182156580918
$ws.Range("G6").Value = 10

# --- Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 30
$ws.Range("F2").Value = 0.9099892365704194
$ws.Range("G2").Value = 19
$ws.Range("D3").Value = 13
$ws.Range("F3").Value = 0.9146017828343822
$ws.Range("G3").Value = 11
$ws.Range("F4").Value = 0.9376483026903331
$ws.Range("D5").Value = 30
$ws.Range("F5").Value = 0.9249282609316706
$ws.Range("G5").Value = 20
$ws.Range("D6").Value = 22
$ws.Range("F6").Value = 0.9107486625985826
$ws.Range("G6").Value = 20
